# act tablas web jul25
# Adds the 2023 and 2022 data points to the "Data" sheet (shifting existing
# rows down) and records the "actualizacion" / "Julio 2025" metadata entry
# on the "Metadata" sheet (shifting the citation rows down).

$wb = $excel.ActiveWorkbook

# --- Data sheet: insert the two newest years at the top of the series ---
$ws1 = $wb.Worksheets.Item("Data")
$ws1.Rows("2:3").Insert()

$ws1.Range("A2").Value = "'2023"
$ws1.Range("B2").Value = 39.7

$ws1.Range("A3").Value = "'2022"
$ws1.Range("B3").Value = 38.3

# --- Metadata sheet: fix the blank placeholder row and add the update note ---
$ws2 = $wb.Worksheets.Item("Metadata")
$ws2.Range("A1").Value = " "

$ws2.Rows("9:9").Insert()
$ws2.Range("A9").Value = "actualizacion"
$ws2.Range("B9").Value = "Julio 2025"
